# CodeSystem-footstrike-type-cs.xlsx — regenerate the "Metadata" sheet's
# Experimental flag and Date stamp (FHIR IG export refresh).
#
# Row 7  "Experimental" -> Value column (B7) should read the literal text
#         "false" (not the Excel boolean FALSE).
# Row 8  "Date"         -> Value column (B8) timestamp bumped to the new
#         generation time.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Metadata")

# Plain `.Value = "false"` gets auto-coerced by Excel into the boolean
# FALSE (same as typing `false` into a cell). We need literal text, so we
# write a text-returning formula and then paste back just the value,
# which downgrades the formula cell into a plain inline/shared string
# while preserving the existing cell style.
$ws.Range("B7").Formula = '=T("false")'
$ws.Range("B7").Copy()
$ws.Range("B7").PasteSpecial(-4163)  # xlPasteValues
$excel.CutCopyMode = $false

# Straightforward text replacement for the generation timestamp.
$ws.Range("B8").Value = "2025-11-30T13:08:37+00:00"
